# dictator patch gui clean up
#
# Append trades 97-101 to "Trade Log" and retype row 96's percentage
# columns (M/N) from text to numeric. Row 101 is the only new row whose
# percentage cells keep the legacy text formatting, so its M/N values are
# copied straight from row 96's original (pre-edit) text cells before
# those get overwritten with numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 101 keeps the old text-typed "0.02" / "0.01" -- grab them from row 96
# while it still holds the original text cells.
$ws.Cells.Item(96, 13).Copy($ws.Cells.Item(101, 13))
$ws.Cells.Item(96, 14).Copy($ws.Cells.Item(101, 14))

# Columns: TradeID, Timestamp, Symbol, TradeType, Price, MarketType, Status,
#          MarketCondition, ma_200, ma_21, ma_7, ma_5, ReturnPct, LossRiskPct,
#          Interval
$rows = @(
    @(96, "a8c9bea8-5b2b-4c2d-abe7-add23ea1ab76", "2025-03-06 14:57:37", "btcUSDT", "market_buy",
      91484.10000000001, "futures", "Closed", "Bullish",
      91050.6875, 91487.01904761905, 91426.82857142857, 91435.64,
      0.02, 0.01, "1m"),
    @(97, "dd25c90d-7a65-4bb3-92c5-3f294af53f96", "2025-03-10 13:33:28", "btcUSDT", "market_buy",
      82327.60000000001, "futures", "Opened", "Bullish",
      81152.6565, 82184.21904761906, 82201.38571428573, 82220.48000000001,
      0.02, 0.01, "1m"),
    @(98, "dd25c90d-7a65-4bb3-92c5-3f294af53f96", "2025-03-10 13:33:36", "btcUSDT", "market_buy",
      82349.2, "futures", "Closed", "Bullish",
      81152.65700000001, 82184.22380952381, 82201.40000000001, 82220.5,
      0.02, 0.01, "1m"),
    @(99, "c2e6ab33-b3c0-47e4-be70-2344fbc4d023", "2025-03-10 13:33:42", "btcUSDT", "market_buy",
      82316.5, "futures", "Opened", "Bullish",
      81152.66, 82184.25238095238, 82201.48571428572, 82220.62,
      0.02, 0.01, "1m"),
    @(100, "c2e6ab33-b3c0-47e4-be70-2344fbc4d023", "2025-03-10 13:33:48", "btcUSDT", "market_buy",
      82357.39999999999, "futures", "Closed", "Bullish",
      81152.811, 82185.69047619047, 82205.80000000002, 82226.66,
      0.02, 0.01, "1m")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowNum = $r[0]

    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
    $ws.Cells.Item($rowNum, 9).Value = $r[9]
    $ws.Cells.Item($rowNum, 10).Value = $r[10]
    $ws.Cells.Item($rowNum, 11).Value = $r[11]
    $ws.Cells.Item($rowNum, 12).Value = $r[12]
    $ws.Cells.Item($rowNum, 13).Value = $r[13]
    $ws.Cells.Item($rowNum, 14).Value = $r[14]
    $ws.Cells.Item($rowNum, 15).Value = $r[15]
}

# Row 101: new trade; M101/N101 were already seeded (text) above.
$ws.Cells.Item(101, 1).Value = "2e905a54-7391-4929-80f3-05f76ca1d71a"
$ws.Cells.Item(101, 2).Value = "2025-03-10 18:55:45"
$ws.Cells.Item(101, 3).Value = "btcUSDT"
$ws.Cells.Item(101, 4).Value = "market_buy"
$ws.Cells.Item(101, 5).Value = 82265.5
$ws.Cells.Item(101, 6).Value = "futures"
$ws.Cells.Item(101, 7).Value = "Opened"
$ws.Cells.Item(101, 8).Value = "Shortbull"
$ws.Cells.Item(101, 9).Value = 82298.94
$ws.Cells.Item(101, 10).Value = 82226.8142857143
$ws.Cells.Item(101, 11).Value = 82256.67142857144
$ws.Cells.Item(101, 12).Value = 82269.56
$ws.Cells.Item(101, 15).Value = "1m"
